$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 77: fill in the export-data columns (E:K) for the first new export file
$ws.Range("E77").Value = "18/05/2023"
$ws.Range("F77").Value = "Blog"
$ws.Range("G77").Value = "5"
$ws.Range("H77").Value = "superblogs.com.ua"
$ws.Range("I77").Value = "https:/superblogs.com.ua/post2.html"
$ws.Range("J77").Value = "таможенне оформления 3"
$ws.Range("K77").Value = 12

# Row 78: fill in the export-data columns (E:K) for the second new export file
$ws.Range("E78").Value = "15/03/2023"
$ws.Range("F78").Value = "Blog"
$ws.Range("G78").Value = "5"
$ws.Range("H78").Value = "superblogs.com.ua"
$ws.Range("I78").Value = "https:/superblogs.com.ua/post3.html"
$ws.Range("J78").Value = "таможенне оформления 4"
$ws.Range("K78").Value = 13

# Row 79: new blank row appended at the end (all cells empty string)
$ws.Range("A79").Value = ""
$ws.Range("B79").Value = ""
$ws.Range("C79").Value = ""
$ws.Range("D79").Value = ""
$ws.Range("E79").Value = ""
$ws.Range("F79").Value = ""
$ws.Range("G79").Value = ""
$ws.Range("H79").Value = ""
$ws.Range("I79").Value = ""
$ws.Range("J79").Value = ""
$ws.Range("K79").Value = ""
